$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 108.818184
$ws.Range("I9").Value = 123.375
$ws.Range("J9").Value = 70
$ws.Range("K9").Value = 123.375
$ws.Range("L9").Value = 70
$ws.Range("M9").Value = 45.625
$ws.Range("N9").Value = -408

$ws.Range("H17").Value = 247376.81
$ws.Range("J17").Value = 252561.22
$ws.Range("L17").Value = 757683.66
$ws.Range("N17").Value = -758019.66

$ws.Range("H33").Value = 10870154
$ws.Range("I33").Value = 277.0625
$ws.Range("K33").Value = 277.0625
$ws.Range("M33").Value = -48.0625

$ws.Range("H76").Value = 21085396
$ws.Range("I76").Value = 3250
$ws.Range("K76").Value = 3250
$ws.Range("M76").Value = -2935

$ws.Range("H79").Value = 21085396
$ws.Range("I79").Value = 3250
$ws.Range("K79").Value = 3250
$ws.Range("M79").Value = -2158

$ws.Range("H92").Value = 112366.336
$ws.Range("I92").Value = 1542.5714
$ws.Range("J92").Value = 500249.5
$ws.Range("K92").Value = 1542.5714
$ws.Range("L92").Value = 500249.5
$ws.Range("M92").Value = -294.5714
$ws.Range("N92").Value = -502745.5

$ws.Range("H137").Value = 2517.2092
$ws.Range("I137").Value = 2093.75
$ws.Range("K137").Value = 6281.25
$ws.Range("M137").Value = -3731.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2015
$ws.Range("I2").Value = 1986.5625
$ws.Range("J2").Value = 2166.6667
$ws.Range("K2").Value = 1986.5625
$ws.Range("L2").Value = 2166.6667
$ws.Range("M2").Value = -1873.5625
$ws.Range("N2").Value = -2392.6667

$ws.Range("H61").Value = 71430616
$ws.Range("I61").Value = 83335224
$ws.Range("K61").Value = 83335224
$ws.Range("M61").Value = -83335012

$ws.Range("H88").Value = 8773226
$ws.Range("I88").Value = 27778182
$ws.Range("K88").Value = 27778182
$ws.Range("M88").Value = -27777776

$ws.Range("H91").Value = 8773226
$ws.Range("I91").Value = 27778182
$ws.Range("K91").Value = 27778182
$ws.Range("M91").Value = -27776778

$ws.Range("H102").Value = 41537.83
$ws.Range("I102").Value = 42536.273
$ws.Range("J102").Value = 38399.855
$ws.Range("K102").Value = 42536.273
$ws.Range("L102").Value = 38399.855
$ws.Range("M102").Value = -40914.273
$ws.Range("N102").Value = -41643.855

$ws.Range("H116").Value = 2015
$ws.Range("I116").Value = 1986.5625
$ws.Range("J116").Value = 2166.6667
$ws.Range("K116").Value = 1986.5625
$ws.Range("L116").Value = 2166.6667
$ws.Range("M116").Value = 307.4375
$ws.Range("N116").Value = -6754.6667

$ws.Range("H132").Value = 2517.6365
$ws.Range("I132").Value = 1589.1428
$ws.Range("J132").Value = 4142.5
$ws.Range("K132").Value = 4767.428400000001
$ws.Range("L132").Value = 12427.5
$ws.Range("M132").Value = -2237.428400000001
$ws.Range("N132").Value = -17487.5

$ws.Range("H136").Value = 71430616
$ws.Range("I136").Value = 83335224
$ws.Range("K136").Value = 250005672
$ws.Range("M136").Value = -250003122

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2015
$ws.Range("I3").Value = 1986.5625
$ws.Range("J3").Value = 2166.6667
$ws.Range("K3").Value = 1986.5625
$ws.Range("L3").Value = 2166.6667
$ws.Range("M3").Value = -1872.5625
$ws.Range("N3").Value = -2394.6667

$ws.Range("H11").Value = 1735.25
$ws.Range("I11").Value = 312
$ws.Range("J11").Value = 6005
$ws.Range("K11").Value = 312
$ws.Range("L11").Value = 6005
$ws.Range("M11").Value = -172
$ws.Range("N11").Value = -6285

$ws.Range("H20").Value = 13953.031
$ws.Range("I20").Value = 12603.84
$ws.Range("K20").Value = 12603.84
$ws.Range("M20").Value = -12356.84

$ws.Range("H80").Value = 33333570
$ws.Range("I80").Value = 354.8
$ws.Range("K80").Value = 354.8
$ws.Range("M80").Value = 643.2

$ws.Range("H83").Value = 33333570
$ws.Range("I83").Value = 354.8
$ws.Range("K83").Value = 1774
$ws.Range("M83").Value = 3218

$ws.Range("H86").Value = 3190.12
$ws.Range("J86").Value = 1631.625
$ws.Range("L86").Value = 1631.625
$ws.Range("N86").Value = -3877.625

$ws.Range("H89").Value = 3190.12
$ws.Range("J89").Value = 1631.625
$ws.Range("L89").Value = 8158.125
$ws.Range("N89").Value = -19390.125

$ws.Range("H107").Value = 55001988
$ws.Range("I107").Value = 7144671.5
$ws.Range("K107").Value = 7144671.5
$ws.Range("M107").Value = -7142751.5

$ws.Range("H134").Value = 1463.8276
$ws.Range("I134").Value = 899.4400000000001
$ws.Range("J134").Value = 4991.25
$ws.Range("K134").Value = 2698.32
$ws.Range("L134").Value = 14973.75
$ws.Range("M134").Value = -163.3200000000002
$ws.Range("N134").Value = -20043.75

$ws.Range("H135").Value = 45265.145
$ws.Range("J135").Value = 45265.145
$ws.Range("L135").Value = 45265.145
$ws.Range("N135").Value = -55405.145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1770.9259
$ws.Range("I58").Value = 1075.0476
$ws.Range("K58").Value = 1075.0476
$ws.Range("M58").Value = -872.0476000000001

$ws.Range("H99").Value = 795750.4399999999
$ws.Range("I99").Value = 927975.7
$ws.Range("K99").Value = 927975.7
$ws.Range("M99").Value = -926477.7

$ws.Range("H107").Value = 3104.0688
$ws.Range("I107").Value = 2218.6667
$ws.Range("J107").Value = 5428.25
$ws.Range("K107").Value = 2218.6667
$ws.Range("L107").Value = 5428.25
$ws.Range("M107").Value = -298.6667000000002
$ws.Range("N107").Value = -9268.25

$ws.Range("H126").Value = 795750.4399999999
$ws.Range("I126").Value = 927975.7
$ws.Range("K126").Value = 2783927.1
$ws.Range("M126").Value = -2781457.1

$ws.Range("H132").Value = 38719.31
$ws.Range("I132").Value = 69664.86
$ws.Range("K132").Value = 208994.58
$ws.Range("M132").Value = -206464.58

$ws.Range("H134").Value = 3242.5
$ws.Range("I134").Value = 1537.3334
$ws.Range("K134").Value = 4612.0002
$ws.Range("M134").Value = -2077.0002

$ws.Range("H136").Value = 1770.9259
$ws.Range("I136").Value = 1075.0476
$ws.Range("K136").Value = 3225.142800000001
$ws.Range("M136").Value = -675.1428000000005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 982.7273
$ws.Range("I97").Value = 900
$ws.Range("J97").Value = 1030
$ws.Range("K97").Value = 2700
$ws.Range("L97").Value = 3090
$ws.Range("M97").Value = -2204
$ws.Range("N97").Value = -4082

$ws.Range("H136").Value = 2315
$ws.Range("I136").Value = 2315
$ws.Range("K136").Value = 6945
$ws.Range("M136").Value = -1845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 118.375
$ws.Range("I2").Value = 92.42856999999999
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 92.42856999999999
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = 20.57143000000001
$ws.Range("N2").Value = -526

$ws.Range("H80").Value = 2937.923
$ws.Range("I80").Value = 3374.375
$ws.Range("J80").Value = 2239.6
$ws.Range("K80").Value = 3374.375
$ws.Range("L80").Value = 2239.6
$ws.Range("M80").Value = -2376.375
$ws.Range("N80").Value = -4235.6

$ws.Range("H83").Value = 2937.923
$ws.Range("I83").Value = 3374.375
$ws.Range("J83").Value = 2239.6
$ws.Range("K83").Value = 16871.875
$ws.Range("L83").Value = 11198
$ws.Range("M83").Value = -11879.875
$ws.Range("N83").Value = -21182

$ws.Range("H97").Value = 1475.7693
$ws.Range("I97").Value = 1266.84
$ws.Range("J97").Value = 1848.8572
$ws.Range("K97").Value = 1266.84
$ws.Range("L97").Value = 1848.8572
$ws.Range("M97").Value = -770.8399999999999
$ws.Range("N97").Value = -2840.8572

$ws.Range("H132").Value = 5838.7334
$ws.Range("I132").Value = 5067.12
$ws.Range("J132").Value = 9696.799999999999
$ws.Range("K132").Value = 15201.36
$ws.Range("L132").Value = 29090.4
$ws.Range("M132").Value = -12671.36
$ws.Range("N132").Value = -34150.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3451.4602
$ws.Range("I132").Value = 1995.4773
$ws.Range("J132").Value = 6823.2104
$ws.Range("K132").Value = 5986.4319
$ws.Range("L132").Value = 20469.6312
$ws.Range("M132").Value = -3456.4319
$ws.Range("N132").Value = -25529.6312

$ws.Range("H133").Value = 24263
$ws.Range("J133").Value = 24263
$ws.Range("L133").Value = 24263
$ws.Range("N133").Value = -29323

$ws.Range("H136").Value = 1832.1666
$ws.Range("I136").Value = 1709.9111
$ws.Range("J136").Value = 3666
$ws.Range("K136").Value = 5129.7333
$ws.Range("L136").Value = 10998
$ws.Range("M136").Value = -2579.7333
$ws.Range("N136").Value = -16098

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 2606.24
$ws.Range("I18").Value = 1272.7273
$ws.Range("J18").Value = 3654
$ws.Range("K18").Value = 1272.7273
$ws.Range("L18").Value = 3654
$ws.Range("M18").Value = -1099.7273
$ws.Range("N18").Value = -4000

$ws.Range("H54").Value = 9999.666999999999
$ws.Range("J54").Value = 9999.666999999999
$ws.Range("L54").Value = 9999.666999999999
$ws.Range("N54").Value = -11039.667

$ws.Range("H98").Value = 93989
$ws.Range("J98").Value = 93989
$ws.Range("L98").Value = 93989
$ws.Range("N98").Value = -99979

$ws.Range("H126").Value = 2251.2
$ws.Range("I126").Value = 772.5
$ws.Range("J126").Value = 2415.5
$ws.Range("K126").Value = 2317.5
$ws.Range("L126").Value = 7246.5
$ws.Range("M126").Value = 152.5
$ws.Range("N126").Value = -12186.5

$ws.Range("H132").Value = 3665.4146
$ws.Range("I132").Value = 2704.4194
$ws.Range("K132").Value = 8113.2582
$ws.Range("M132").Value = -5583.2582

$ws.Range("H136").Value = 2678.1943
$ws.Range("I136").Value = 1853.7931
$ws.Range("J136").Value = 6093.5713
$ws.Range("K136").Value = 5561.379300000001
$ws.Range("L136").Value = 18280.7139
$ws.Range("M136").Value = -3011.379300000001
$ws.Range("N136").Value = -23380.7139
